# Update B1 header text casing: "Lotsize" -> "LotSize"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "LotSize"

# Add new column F with "Expiry" header and a text-formatted date value below it
$ws.Range("F1").Value = "Expiry"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2").Value = "24JUL2024"

# Center (horizontally & vertically) all the header and data cells in the used range
$usedRange = $ws.Range("A1:F2")
$usedRange.HorizontalAlignment = -4108
$usedRange.VerticalAlignment = -4108

# Store the expiry date as text (keeps the "24JUL2024" literal instead of a date serial)
$ws.Range("F2").NumberFormat = "@"

# Widen column F to fit the new "Expiry" header/value
$ws.Columns.Item(6).ColumnWidth = 19.44140625

# Move the active selection to B2, matching the saved view state
$ws.Range("B2").Select()
